$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Cells.Item(2, 8).Value = 1125.7142
$ws.Cells.Item(2, 9).Value = 1125.7142
$ws.Cells.Item(2, 10).Value = 0
$ws.Cells.Item(2, 11).Value = 1125.7142
$ws.Cells.Item(2, 12).Value = 0
$ws.Cells.Item(2, 13).Value = -1012.7142
$ws.Cells.Item(2, 14).ClearContents()
# Row 32
$ws.Cells.Item(32, 8).Value = 9166
$ws.Cells.Item(32, 9).Value = 9498
$ws.Cells.Item(32, 10).Value = 9000
$ws.Cells.Item(32, 11).Value = 9498
$ws.Cells.Item(32, 12).Value = 9000
$ws.Cells.Item(32, 13).Value = -9172
$ws.Cells.Item(32, 14).Value = -9652
# Row 64
$ws.Cells.Item(64, 8).Value = 4219.115
$ws.Cells.Item(64, 9).Value = 4009
$ws.Cells.Item(64, 10).Value = 5374.75
$ws.Cells.Item(64, 11).Value = 4009
$ws.Cells.Item(64, 12).Value = 5374.75
$ws.Cells.Item(64, 13).Value = -3761
$ws.Cells.Item(64, 14).Value = -5870.75
# Row 67
$ws.Cells.Item(67, 8).Value = 4219.115
$ws.Cells.Item(67, 9).Value = 4009
$ws.Cells.Item(67, 10).Value = 5374.75
$ws.Cells.Item(67, 11).Value = 4009
$ws.Cells.Item(67, 12).Value = 5374.75
$ws.Cells.Item(67, 13).Value = -3151
$ws.Cells.Item(67, 14).Value = -7090.75
# Row 86
$ws.Cells.Item(86, 8).Value = 3258.2
$ws.Cells.Item(86, 9).Value = 3130.3333
$ws.Cells.Item(86, 10).Value = 3450
$ws.Cells.Item(86, 11).Value = 3130.3333
$ws.Cells.Item(86, 12).Value = 3450
$ws.Cells.Item(86, 13).Value = -2007.3333
$ws.Cells.Item(86, 14).Value = -5696
# Row 87
$ws.Cells.Item(87, 8).Value = 33333.75
$ws.Cells.Item(87, 9).Value = 0
$ws.Cells.Item(87, 10).Value = 33333.75
$ws.Cells.Item(87, 11).Value = 0
$ws.Cells.Item(87, 12).Value = 33333.75
$ws.Cells.Item(87, 14).Value = -35829.75
# Row 89
$ws.Cells.Item(89, 8).Value = 3258.2
$ws.Cells.Item(89, 9).Value = 3130.3333
$ws.Cells.Item(89, 10).Value = 3450
$ws.Cells.Item(89, 11).Value = 15651.6665
$ws.Cells.Item(89, 12).Value = 17250
$ws.Cells.Item(89, 13).Value = -10035.6665
$ws.Cells.Item(89, 14).Value = -28482
# Row 90
$ws.Cells.Item(90, 8).Value = 33333.75
$ws.Cells.Item(90, 9).Value = 0
$ws.Cells.Item(90, 10).Value = 33333.75
$ws.Cells.Item(90, 11).Value = 0
$ws.Cells.Item(90, 12).Value = 100001.25
$ws.Cells.Item(90, 14).Value = -112481.25
# Row 98
$ws.Cells.Item(98, 8).Value = 837.45
$ws.Cells.Item(98, 9).Value = 837.45
$ws.Cells.Item(98, 10).Value = 0
$ws.Cells.Item(98, 11).Value = 837.45
$ws.Cells.Item(98, 12).Value = 0
$ws.Cells.Item(98, 13).Value = 660.55
$ws.Cells.Item(98, 14).ClearContents()
# Row 122
$ws.Cells.Item(122, 8).Value = 837.45
$ws.Cells.Item(122, 9).Value = 837.45
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 11).Value = 2512.35
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 13).Value = -62.35000000000036
$ws.Cells.Item(122, 14).ClearContents()
# Row 132
$ws.Cells.Item(132, 8).Value = 1212.7059
$ws.Cells.Item(132, 9).Value = 828.70215
$ws.Cells.Item(132, 10).Value = 5724.75
$ws.Cells.Item(132, 11).Value = 2486.10645
$ws.Cells.Item(132, 12).Value = 17174.25
$ws.Cells.Item(132, 13).Value = 43.89355000000023
$ws.Cells.Item(132, 14).Value = -22234.25
# Row 137
$ws.Cells.Item(137, 8).Value = 2181.0645
$ws.Cells.Item(137, 9).Value = 2126.3333
$ws.Cells.Item(137, 10).Value = 2368.7144
$ws.Cells.Item(137, 11).Value = 6378.999899999999
$ws.Cells.Item(137, 12).Value = 7106.1432
$ws.Cells.Item(137, 13).Value = -3828.999899999999
$ws.Cells.Item(137, 14).Value = -12206.1432
# Row 138
$ws.Cells.Item(138, 8).Value = 2226.2341
$ws.Cells.Item(138, 9).Value = 1475.1428
$ws.Cells.Item(138, 10).Value = 2832.8845
$ws.Cells.Item(138, 11).Value = 4425.428400000001
$ws.Cells.Item(138, 12).Value = 8498.6535
$ws.Cells.Item(138, 13).Value = 714.5715999999993
$ws.Cells.Item(138, 14).Value = -18778.6535

$ws = $wb.Worksheets.Item("ARM")
# Row 41
$ws.Cells.Item(41, 8).Value = 8098.8
$ws.Cells.Item(41, 9).Value = 2623.5
$ws.Cells.Item(41, 10).Value = 30000
$ws.Cells.Item(41, 11).Value = 2623.5
$ws.Cells.Item(41, 12).Value = 30000
$ws.Cells.Item(41, 13).Value = -2209.5
$ws.Cells.Item(41, 14).Value = -30828
# Row 132
$ws.Cells.Item(132, 8).Value = 2144.0571
$ws.Cells.Item(132, 9).Value = 1140.8966
$ws.Cells.Item(132, 10).Value = 6992.6665
$ws.Cells.Item(132, 11).Value = 3422.6898
$ws.Cells.Item(132, 12).Value = 20977.9995
$ws.Cells.Item(132, 13).Value = -892.6898000000001
$ws.Cells.Item(132, 14).Value = -26037.9995

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Cells.Item(20, 8).Value = 6706.07
$ws.Cells.Item(20, 9).Value = 5767.654
$ws.Cells.Item(20, 10).Value = 8141.294
$ws.Cells.Item(20, 11).Value = 5767.654
$ws.Cells.Item(20, 12).Value = 8141.294
$ws.Cells.Item(20, 13).Value = -5520.654
$ws.Cells.Item(20, 14).Value = -8635.294
# Row 86
$ws.Cells.Item(86, 8).Value = 2908.8215
$ws.Cells.Item(86, 9).Value = 2104.5557
$ws.Cells.Item(86, 10).Value = 3289.7896
$ws.Cells.Item(86, 11).Value = 2104.5557
$ws.Cells.Item(86, 12).Value = 3289.7896
$ws.Cells.Item(86, 13).Value = -981.5556999999999
$ws.Cells.Item(86, 14).Value = -5535.7896
# Row 89
$ws.Cells.Item(89, 8).Value = 2908.8215
$ws.Cells.Item(89, 9).Value = 2104.5557
$ws.Cells.Item(89, 10).Value = 3289.7896
$ws.Cells.Item(89, 11).Value = 10522.7785
$ws.Cells.Item(89, 12).Value = 16448.948
$ws.Cells.Item(89, 13).Value = -4906.7785
$ws.Cells.Item(89, 14).Value = -27680.948
# Row 94
$ws.Cells.Item(94, 8).Value = 1023.63635
$ws.Cells.Item(94, 9).Value = 1023.63635
$ws.Cells.Item(94, 10).Value = 0
$ws.Cells.Item(94, 11).Value = 1023.63635
$ws.Cells.Item(94, 12).Value = 0
$ws.Cells.Item(94, 13).Value = -572.63635
# Row 99
$ws.Cells.Item(99, 8).Value = 1699.6428
$ws.Cells.Item(99, 9).Value = 1407.3077
$ws.Cells.Item(99, 10).Value = 5500
$ws.Cells.Item(99, 11).Value = 1407.3077
$ws.Cells.Item(99, 12).Value = 5500
$ws.Cells.Item(99, 13).Value = 90.69229999999993
$ws.Cells.Item(99, 14).Value = -8496
# Row 105
$ws.Cells.Item(105, 8).Value = 9578.777
$ws.Cells.Item(105, 9).Value = 10429.857
$ws.Cells.Item(105, 10).Value = 6600
$ws.Cells.Item(105, 11).Value = 10429.857
$ws.Cells.Item(105, 12).Value = 6600
$ws.Cells.Item(105, 13).Value = -8682.857
$ws.Cells.Item(105, 14).Value = -10094
# Row 107
$ws.Cells.Item(107, 8).Value = 1708.8823
$ws.Cells.Item(107, 9).Value = 1860.7858
$ws.Cells.Item(107, 10).Value = 1000
$ws.Cells.Item(107, 11).Value = 1860.7858
$ws.Cells.Item(107, 12).Value = 1000
$ws.Cells.Item(107, 13).Value = 59.21419999999989
$ws.Cells.Item(107, 14).Value = -4840

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Cells.Item(31, 8).Value = 2200.35
$ws.Cells.Item(31, 9).Value = 1499.6666
$ws.Cells.Item(31, 10).Value = 8506.5
$ws.Cells.Item(31, 11).Value = 1499.6666
$ws.Cells.Item(31, 12).Value = 8506.5
$ws.Cells.Item(31, 13).Value = -1204.6666
$ws.Cells.Item(31, 14).Value = -9096.5
# Row 34
$ws.Cells.Item(34, 8).Value = 2200.35
$ws.Cells.Item(34, 9).Value = 1499.6666
$ws.Cells.Item(34, 10).Value = 8506.5
$ws.Cells.Item(34, 11).Value = 1499.6666
$ws.Cells.Item(34, 12).Value = 8506.5
$ws.Cells.Item(34, 13).Value = -1297.6666
$ws.Cells.Item(34, 14).Value = -8910.5
# Row 132
$ws.Cells.Item(132, 8).Value = 1762.4706
$ws.Cells.Item(132, 9).Value = 1710.5161
$ws.Cells.Item(132, 10).Value = 2299.3333
$ws.Cells.Item(132, 11).Value = 5131.5483
$ws.Cells.Item(132, 12).Value = 6897.999899999999
$ws.Cells.Item(132, 13).Value = -2601.5483
$ws.Cells.Item(132, 14).Value = -11957.9999

$ws = $wb.Worksheets.Item("CUL")
# Row 17
$ws.Cells.Item(17, 8).Value = 125.19231
$ws.Cells.Item(17, 9).Value = 125.19231
$ws.Cells.Item(17, 10).Value = 0
$ws.Cells.Item(17, 11).Value = 375.57693
$ws.Cells.Item(17, 12).Value = 0
$ws.Cells.Item(17, 13).Value = -206.57693
# Row 34
$ws.Cells.Item(34, 8).Value = 1417.8948
$ws.Cells.Item(34, 9).Value = 115.833336
$ws.Cells.Item(34, 10).Value = 2018.8462
$ws.Cells.Item(34, 11).Value = 347.500008
$ws.Cells.Item(34, 12).Value = 6056.5386
$ws.Cells.Item(34, 13).Value = -263.500008
$ws.Cells.Item(34, 14).Value = -6224.5386
# Row 40
$ws.Cells.Item(40, 8).Value = 68.40000000000001
$ws.Cells.Item(40, 9).Value = 57.857143
$ws.Cells.Item(40, 10).Value = 93
$ws.Cells.Item(40, 11).Value = 231.428572
$ws.Cells.Item(40, 12).Value = 372
$ws.Cells.Item(40, 13).Value = -162.428572
$ws.Cells.Item(40, 14).Value = -510
# Row 55
$ws.Cells.Item(55, 8).Value = 11365605
$ws.Cells.Item(55, 9).Value = 615.5
$ws.Cells.Item(55, 10).Value = 17859886
$ws.Cells.Item(55, 11).Value = 1846.5
$ws.Cells.Item(55, 12).Value = 53579658
$ws.Cells.Item(55, 13).Value = -1669.5
$ws.Cells.Item(55, 14).Value = -53580012
# Row 69
$ws.Cells.Item(69, 8).Value = 2999.5
$ws.Cells.Item(69, 9).Value = 0
$ws.Cells.Item(69, 10).Value = 2999.5
$ws.Cells.Item(69, 11).Value = 0
$ws.Cells.Item(69, 12).Value = 8998.5
$ws.Cells.Item(69, 14).Value = -10620.5
# Row 72
$ws.Cells.Item(72, 8).Value = 2999.5
$ws.Cells.Item(72, 9).Value = 0
$ws.Cells.Item(72, 10).Value = 2999.5
$ws.Cells.Item(72, 11).Value = 0
$ws.Cells.Item(72, 12).Value = 26995.5
$ws.Cells.Item(72, 14).Value = -35107.5
# Row 128
$ws.Cells.Item(128, 8).Value = 443339
$ws.Cells.Item(128, 9).Value = 443339
$ws.Cells.Item(128, 10).Value = 0
$ws.Cells.Item(128, 11).Value = 1330017
$ws.Cells.Item(128, 12).Value = 0
$ws.Cells.Item(128, 13).Value = -1325037
# Row 140
$ws.Cells.Item(140, 8).Value = 3173
$ws.Cells.Item(140, 9).Value = 2361.3333
$ws.Cells.Item(140, 10).Value = 4999.25
$ws.Cells.Item(140, 11).Value = 7083.999899999999
$ws.Cells.Item(140, 12).Value = 14997.75
$ws.Cells.Item(140, 13).Value = -1903.999899999999
$ws.Cells.Item(140, 14).Value = -25357.75

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Cells.Item(70, 8).Value = 4351.143
$ws.Cells.Item(70, 9).Value = 3287
$ws.Cells.Item(70, 10).Value = 5415.2856
$ws.Cells.Item(70, 11).Value = 3287
$ws.Cells.Item(70, 12).Value = 5415.2856
$ws.Cells.Item(70, 13).Value = -3017
$ws.Cells.Item(70, 14).Value = -5955.2856
# Row 73
$ws.Cells.Item(73, 8).Value = 4351.143
$ws.Cells.Item(73, 9).Value = 3287
$ws.Cells.Item(73, 10).Value = 5415.2856
$ws.Cells.Item(73, 11).Value = 3287
$ws.Cells.Item(73, 12).Value = 5415.2856
$ws.Cells.Item(73, 13).Value = -2351
$ws.Cells.Item(73, 14).Value = -7287.2856
# Row 97
$ws.Cells.Item(97, 8).Value = 36190.95
$ws.Cells.Item(97, 9).Value = 60643.816
$ws.Cells.Item(97, 10).Value = 2568.25
$ws.Cells.Item(97, 11).Value = 60643.816
$ws.Cells.Item(97, 12).Value = 2568.25
$ws.Cells.Item(97, 13).Value = -60147.816
$ws.Cells.Item(97, 14).Value = -3560.25
# Row 122
$ws.Cells.Item(122, 8).Value = 3007.7222
$ws.Cells.Item(122, 9).Value = 3093
$ws.Cells.Item(122, 10).Value = 2581.3333
$ws.Cells.Item(122, 11).Value = 9279
$ws.Cells.Item(122, 12).Value = 7743.999899999999
$ws.Cells.Item(122, 13).Value = -6829
$ws.Cells.Item(122, 14).Value = -12643.9999
# Row 132
$ws.Cells.Item(132, 8).Value = 2246.2173
$ws.Cells.Item(132, 9).Value = 2034.65
$ws.Cells.Item(132, 10).Value = 3656.6667
$ws.Cells.Item(132, 11).Value = 6103.950000000001
$ws.Cells.Item(132, 12).Value = 10970.0001
$ws.Cells.Item(132, 13).Value = -3573.950000000001
$ws.Cells.Item(132, 14).Value = -16030.0001

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Cells.Item(7, 8).Value = 32929.332
$ws.Cells.Item(7, 9).Value = 46398.168
$ws.Cells.Item(7, 10).Value = 5991.6665
$ws.Cells.Item(7, 11).Value = 46398.168
$ws.Cells.Item(7, 12).Value = 5991.6665
$ws.Cells.Item(7, 13).Value = -46286.168
$ws.Cells.Item(7, 14).Value = -6215.6665
# Row 55
$ws.Cells.Item(55, 8).Value = 673.05
$ws.Cells.Item(55, 9).Value = 554
$ws.Cells.Item(55, 10).Value = 818.55554
$ws.Cells.Item(55, 11).Value = 554
$ws.Cells.Item(55, 12).Value = 818.55554
$ws.Cells.Item(55, 13).Value = -381
$ws.Cells.Item(55, 14).Value = -1164.55554
# Row 61
$ws.Cells.Item(61, 8).Value = 16797.238
$ws.Cells.Item(61, 9).Value = 895.55554
$ws.Cells.Item(61, 10).Value = 112207.336
$ws.Cells.Item(61, 11).Value = 895.55554
$ws.Cells.Item(61, 12).Value = 112207.336
$ws.Cells.Item(61, 13).Value = -693.55554
$ws.Cells.Item(61, 14).Value = -112611.336
# Row 113
$ws.Cells.Item(113, 8).Value = 16797.238
$ws.Cells.Item(113, 9).Value = 895.55554
$ws.Cells.Item(113, 10).Value = 112207.336
$ws.Cells.Item(113, 11).Value = 895.55554
$ws.Cells.Item(113, 12).Value = 112207.336
$ws.Cells.Item(113, 13).Value = 1274.44446
$ws.Cells.Item(113, 14).Value = -116547.336
# Row 122
$ws.Cells.Item(122, 8).Value = 2965.2222
$ws.Cells.Item(122, 9).Value = 2752.5417
$ws.Cells.Item(122, 10).Value = 4666.6665
$ws.Cells.Item(122, 11).Value = 8257.625100000001
$ws.Cells.Item(122, 12).Value = 13999.9995
$ws.Cells.Item(122, 13).Value = -5807.625100000001
$ws.Cells.Item(122, 14).Value = -18899.9995
# Row 126
$ws.Cells.Item(126, 8).Value = 32929.332
$ws.Cells.Item(126, 9).Value = 46398.168
$ws.Cells.Item(126, 10).Value = 5991.6665
$ws.Cells.Item(126, 11).Value = 139194.504
$ws.Cells.Item(126, 12).Value = 17974.9995
$ws.Cells.Item(126, 13).Value = -136724.504
$ws.Cells.Item(126, 14).Value = -22914.9995

$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Cells.Item(113, 8).Value = 806.95654
$ws.Cells.Item(113, 9).Value = 378.73334
$ws.Cells.Item(113, 10).Value = 1609.875
$ws.Cells.Item(113, 11).Value = 1136.20002
$ws.Cells.Item(113, 12).Value = 4829.625
$ws.Cells.Item(113, 13).Value = 1033.79998
$ws.Cells.Item(113, 14).Value = -9169.625
# Row 126
$ws.Cells.Item(126, 8).Value = 3395.2632
$ws.Cells.Item(126, 9).Value = 2912.647
$ws.Cells.Item(126, 10).Value = 7497.5
$ws.Cells.Item(126, 11).Value = 8737.940999999999
$ws.Cells.Item(126, 12).Value = 22492.5
$ws.Cells.Item(126, 13).Value = -6267.940999999999
$ws.Cells.Item(126, 14).Value = -27432.5
